# MallaCurricular.xlsx edit
# - Column D ("Prerrequisito") is repurposed into a new column
#   "Abre la/s asignatura/s:" that lists which row IDs a given course
#   unlocks. Some cells hold a single numeric ID, some hold a
#   comma-separated text list of IDs, and some are left blank.
# - Column D keeps its existing per-row fill/border formatting, but is
#   switched to a Text ("@") number format.
# - The sheet view selection moves from D55 to H18, with no frozen
#   top-left-cell scroll offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header (row 2) ----
$ws.Range("D2").Value = "Abre la/s asignatura/s:"

# ---- Column D body values, row by row ----
# Numbers are written as real numeric values (Value, set BEFORE the
# number format is switched to Text) so they stay <v>N</v> cells (no
# t="s") even though the cell's display format is "@".
# Comma lists are written as text (Value set AFTER the number format
# switch), landing as shared-string cells.
# Blank cells just get their contents cleared.

function Set-NumericTextCell($addr, $num) {
    $ws.Range($addr).Value = $num
    $ws.Range($addr).NumberFormat = "@"
}

function Set-TextCell($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

function Clear-Cell($addr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).ClearContents()
}

# Semester 1 (rows 3-7)
Set-NumericTextCell "D3" 6
Set-TextCell         "D4" "7, 8, 33"
Clear-Cell           "D5"
Set-NumericTextCell "D6" 9
Set-NumericTextCell "D7" 38

# Semester 2 (rows 8-12)
Set-TextCell         "D8"  "11, 17, 22"
Set-TextCell         "D9"  "11, 12, 13, 16, 17"
Set-NumericTextCell "D10" 13
Set-TextCell         "D11" "14, 15"
Clear-Cell           "D12"

# Semester 3 (rows 13-17)
Set-NumericTextCell "D13" 18
Set-TextCell         "D14" "18, 22, 28"
Clear-Cell           "D15"
Set-NumericTextCell "D16" 19
Set-TextCell         "D17" "24, 29, 31"

# Semester 4 (rows 18-23)
Set-TextCell         "D18" "29, 30"
Clear-Cell           "D19"
Set-NumericTextCell "D20" 23
Set-TextCell         "D21" "24, 25"
Clear-Cell           "D22"
Set-NumericTextCell "D23" 27

# Semester 5 (rows 24-29)
Set-NumericTextCell "D24" 29
Set-TextCell         "D25" "30, 32"
Set-TextCell         "D26" "34, 41"
Set-TextCell         "D27" "31, 35"
Clear-Cell           "D28"
Clear-Cell           "D29"

# Semester 6 (rows 30-34)
Clear-Cell           "D30"
Clear-Cell           "D31"
Set-NumericTextCell "D32" 35
Set-NumericTextCell "D33" 40
Set-NumericTextCell "D34" 37

# Semester 7 (rows 35-39)
Set-NumericTextCell "D35" 46
Set-TextCell         "D36" "39, 41"
Set-NumericTextCell "D37" 41
Clear-Cell           "D38"
Clear-Cell           "D39"

# Semester 8 (rows 40-44)
Clear-Cell           "D40"
Clear-Cell           "D41"
Clear-Cell           "D42"
Set-NumericTextCell "D43" 46
Clear-Cell           "D44"

# Semester 9 (rows 45-49)
Clear-Cell           "D45"
Clear-Cell           "D46"
Clear-Cell           "D47"
Set-NumericTextCell "D48" 51
Clear-Cell           "D49"

# Semester 10 (rows 50-54)
Clear-Cell           "D50"
Clear-Cell           "D51"
Clear-Cell           "D52"
Clear-Cell           "D53"
Clear-Cell           "D54"

# ---- Sheet view: scroll back to top, select H18 ----
$ws.Range("H18").Select()

Write-Host "Malla curricular updated"
